# Replace the five trailing question cells (A6:A10) with the new set of
# questions. Excel/the engine dedupes the shared-string table on save, so
# re-using "what a belief that divides the world by 50 - 50?" for A6 simply
# re-points it at the existing shared string instead of duplicating it, and
# the four now-orphaned old strings (A6:A9's previous text) are dropped from
# the table automatically, leaving the table in the exact order the new
# cells are written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value  = "what a belief that divides the world by 50 - 50?"
$ws.Range("A7").Value  = "What's the best Instgram comment you read, Explain?"
$ws.Range("A8").Value  = "What's your ""It always happen to good people"" moment?"
$ws.Range("A9").Value  = "Is guilt over ex good or bad?"
$ws.Range("A10").Value = "What's your secret Eureka Moment?"
